# "Add files via upload" — re-upload of the MSFT Models workbook with the
# Historicals tab renamed/expanded and the DCF Model tab left active.

$wb = $excel.ActiveWorkbook

$wsHist = $wb.Worksheets.Item("Historicals")
$wsDcf  = $wb.Worksheets.Item("DCF Model")

# --- Rename "Historicals" -> "Historicals and 3SM" -------------------------
# (Excel automatically repoints every formula that referenced the old name,
# e.g. DCF Model!G19 "=Historicals!G17" -> "='Historicals and 3SM'!G17")
$wsHist.Name = "Historicals and 3SM"

# --- Historicals and 3SM!D5: literal date -> =TODAY() ----------------------
$wsHist.Range("D5").Formula = "=TODAY()"

# --- DCF Model sheet edits --------------------------------------------------
# Row 36 (year header row): the old fill "=F36+1" dragged across G:S is
# replaced starting at H36 with a direct link to the year row above (H25),
# then refilled right.
$wsDcf.Range("G36").ClearContents()
$wsDcf.Range("H36").Formula = "=H25"

# Row 37 (Revenue actuals pulled from row 19): column G's link is cleared.
$wsDcf.Range("G37").ClearContents()

# Row 43 (COGS/EBIT-ish actuals pulled from row 22): column G cleared.
$wsDcf.Range("G43").ClearContents()

# Row 54 (actuals pulled from row 26): column G cleared (H54 keeps its
# pre-existing #REF! error untouched).
$wsDcf.Range("G54").ClearContents()

# Row 55 (% of row 37, driven off row 54): stale H55 (#REF! based) cleared.
$wsDcf.Range("H55").ClearContents()

# Row 57 (actuals pulled from row 29): column G cleared.
$wsDcf.Range("G57").ClearContents()

# Row 58 (% of row 37, driven off row 57): stale H58 cleared.
$wsDcf.Range("H58").ClearContents()

# Row 60 (actuals pulled from row 32): column G cleared.
$wsDcf.Range("G60").ClearContents()

# Row 61 (% of row 37, driven off row 60): stale H61 cleared.
$wsDcf.Range("H61").ClearContents()

# --- Window / selection state -----------------------------------------------
# Historicals and 3SM is no longer the tab in front; DCF Model is now active
# with H1 selected, while Historicals and 3SM's own selection moves to H28.
$wsHist.Range("H28").Select() | Out-Null
$wsDcf.Activate() | Out-Null
$wsDcf.Range("H1").Select() | Out-Null
